$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old data rows (rows 2-7) before writing the new 9-row table (rows 2-10)
$ws.Range("A2:T10").ClearContents()

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Vtn"
$ws.Cells.Item(2,3).Value = "Itga8"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 7.134618
$ws.Cells.Item(2,8).Value = 21.403854
$ws.Cells.Item(2,9).Value = 0.0965317920926077
$ws.Cells.Item(2,10).Value = 0.0965317920926077
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 1.518534666666667
$ws.Cells.Item(2,14).Value = 4.555604
$ws.Cells.Item(2,15).Value = 0.1025715407499064
$ws.Cells.Item(2,16).Value = 0.1025715407499064
$ws.Cells.Item(2,17).Value = 10.834164766424
$ws.Cells.Item(2,18).Value = 97.50748289781599
$ws.Cells.Item(2,19).Value = 0.009901414646288402
$ws.Cells.Item(2,20).Value = 0.009901414646288402

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Vtn"
$ws.Cells.Item(3,3).Value = "Itga8"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 7.134618
$ws.Cells.Item(3,8).Value = 21.403854
$ws.Cells.Item(3,9).Value = 0.0965317920926077
$ws.Cells.Item(3,10).Value = 0.0965317920926077
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 4.418558666666667
$ws.Cells.Item(3,14).Value = 13.255676
$ws.Cells.Item(3,15).Value = 0.2984577041818288
$ws.Cells.Item(3,16).Value = 0.2984577041818288
$ws.Cells.Item(3,17).Value = 31.524728197256
$ws.Cells.Item(3,18).Value = 283.722553775304
$ws.Cells.Item(3,19).Value = 0.02881065704851731
$ws.Cells.Item(3,20).Value = 0.02881065704851731

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Vtn"
$ws.Cells.Item(4,3).Value = "Itga8"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 7.134618
$ws.Cells.Item(4,8).Value = 21.403854
$ws.Cells.Item(4,9).Value = 0.0965317920926077
$ws.Cells.Item(4,10).Value = 0.0965317920926077
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 8.867545999999999
$ws.Cells.Item(4,14).Value = 26.602638
$ws.Cells.Item(4,15).Value = 0.5989707550682648
$ws.Cells.Item(4,16).Value = 0.5989707550682649
$ws.Cells.Item(4,17).Value = 63.26655330742799
$ws.Cells.Item(4,18).Value = 569.398979766852
$ws.Cells.Item(4,19).Value = 0.05781972039780199
$ws.Cells.Item(4,20).Value = 0.057819720397802

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Vtn"
$ws.Cells.Item(5,3).Value = "Itga8"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 17.50798033333334
$ws.Cells.Item(5,8).Value = 52.52394100000001
$ws.Cells.Item(5,9).Value = 0.2368839813846793
$ws.Cells.Item(5,10).Value = 0.2368839813846794
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 1.518534666666667
$ws.Cells.Item(5,14).Value = 4.555604
$ws.Cells.Item(5,15).Value = 0.1025715407499064
$ws.Cells.Item(5,16).Value = 0.1025715407499064
$ws.Cells.Item(5,17).Value = 26.58647507948489
$ws.Cells.Item(5,18).Value = 239.278275715364
$ws.Cells.Item(5,19).Value = 0.0242975549495987
$ws.Cells.Item(5,20).Value = 0.0242975549495987

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Vtn"
$ws.Cells.Item(6,3).Value = "Itga8"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 17.50798033333334
$ws.Cells.Item(6,8).Value = 52.52394100000001
$ws.Cells.Item(6,9).Value = 0.2368839813846793
$ws.Cells.Item(6,10).Value = 0.2368839813846794
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 4.418558666666667
$ws.Cells.Item(6,14).Value = 13.255676
$ws.Cells.Item(6,15).Value = 0.2984577041818288
$ws.Cells.Item(6,16).Value = 0.2984577041818288
$ws.Cells.Item(6,17).Value = 77.36003823767957
$ws.Cells.Item(6,18).Value = 696.240344139116
$ws.Cells.Item(6,19).Value = 0.07069984924152246
$ws.Cells.Item(6,20).Value = 0.07069984924152246

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Vtn"
$ws.Cells.Item(7,3).Value = "Itga8"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 17.50798033333334
$ws.Cells.Item(7,8).Value = 52.52394100000001
$ws.Cells.Item(7,9).Value = 0.2368839813846793
$ws.Cells.Item(7,10).Value = 0.2368839813846794
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 8.867545999999999
$ws.Cells.Item(7,14).Value = 26.602638
$ws.Cells.Item(7,15).Value = 0.5989707550682648
$ws.Cells.Item(7,16).Value = 0.5989707550682649
$ws.Cells.Item(7,17).Value = 155.2528209729287
$ws.Cells.Item(7,18).Value = 1397.275388756358
$ws.Cells.Item(7,19).Value = 0.1418865771935582
$ws.Cells.Item(7,20).Value = 0.1418865771935582

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Vtn"
$ws.Cells.Item(8,3).Value = "Itga8"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 49.26691733333334
$ws.Cells.Item(8,8).Value = 147.800752
$ws.Cells.Item(8,9).Value = 0.6665842265227129
$ws.Cells.Item(8,10).Value = 0.666584226522713
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 1.518534666666667
$ws.Cells.Item(8,14).Value = 4.555604
$ws.Cells.Item(8,15).Value = 0.1025715407499064
$ws.Cells.Item(8,16).Value = 0.1025715407499064
$ws.Cells.Item(8,17).Value = 74.81352189046756
$ws.Cells.Item(8,18).Value = 673.3216970142081
$ws.Cells.Item(8,19).Value = 0.06837257115401926
$ws.Cells.Item(8,20).Value = 0.06837257115401928

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Vtn"
$ws.Cells.Item(9,3).Value = "Itga8"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 49.26691733333334
$ws.Cells.Item(9,8).Value = 147.800752
$ws.Cells.Item(9,9).Value = 0.6665842265227129
$ws.Cells.Item(9,10).Value = 0.666584226522713
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 4.418558666666667
$ws.Cells.Item(9,14).Value = 13.255676
$ws.Cells.Item(9,15).Value = 0.2984577041818288
$ws.Cells.Item(9,16).Value = 0.2984577041818288
$ws.Cells.Item(9,17).Value = 217.6887645631502
$ws.Cells.Item(9,18).Value = 1959.198881068352
$ws.Cells.Item(9,19).Value = 0.198947197891789
$ws.Cells.Item(9,20).Value = 0.198947197891789

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Vtn"
$ws.Cells.Item(10,3).Value = "Itga8"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 49.26691733333334
$ws.Cells.Item(10,8).Value = 147.800752
$ws.Cells.Item(10,9).Value = 0.6665842265227129
$ws.Cells.Item(10,10).Value = 0.666584226522713
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 8.867545999999999
$ws.Cells.Item(10,14).Value = 26.602638
$ws.Cells.Item(10,15).Value = 0.5989707550682648
$ws.Cells.Item(10,16).Value = 0.5989707550682649
$ws.Cells.Item(10,17).Value = 436.8766557315307
$ws.Cells.Item(10,18).Value = 3931.889901583776
$ws.Cells.Item(10,19).Value = 0.3992644574769046
$ws.Cells.Item(10,20).Value = 0.3992644574769048

Write-Host "Updated Vtn-Itga8 sheet with ECs/FAPs/sCs sending x target cluster combinations"